$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 186, shifting existing rows 186:298 down to 187:299
$ws.Rows.Item(186).Insert(-4121)

# Populate the newly inserted row 186 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,Q,R repeat the same market/category metadata
# as the surrounding rows; D,J,K,L,M,O,P are the new data point.
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44603
$ws.Cells.Item(186, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 100114013
$ws.Cells.Item(186, 7).Value = "Zanahoria"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 700
$ws.Cells.Item(186, 11).Value = 11000
$ws.Cells.Item(186, 12).Value = 12000
$ws.Cells.Item(186, 13).Value = 11500
$ws.Cells.Item(186, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(186, 15).Value = "Región de Ñuble"
$ws.Cells.Item(186, 16).Value = 575
$ws.Cells.Item(186, 17).Value = 20
$ws.Cells.Item(186, 18).Value = "Hortaliza"
